$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Inhba"
$ws.Cells.Item(2, 3).Value = "Acvr2a"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 3.675031333333333
$ws.Cells.Item(2, 8).Value = 11.025094
$ws.Cells.Item(2, 9).Value = 0.2032371147293133
$ws.Cells.Item(2, 10).Value = 0.2032371147293133
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 7.940402333333334
$ws.Cells.Item(2, 14).Value = 23.821207
$ws.Cells.Item(2, 15).Value = 0.1931648990487216
$ws.Cells.Item(2, 16).Value = 0.1931648990487216
$ws.Cells.Item(2, 17).Value = 29.18122737427311
$ws.Cells.Item(2, 18).Value = 262.631046368458
$ws.Cells.Item(2, 19).Value = 0.03925827674964126
$ws.Cells.Item(2, 20).Value = 0.03925827674964126

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Inhba"
$ws.Cells.Item(3, 3).Value = "Acvr2a"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 3.675031333333333
$ws.Cells.Item(3, 8).Value = 11.025094
$ws.Cells.Item(3, 9).Value = 0.2032371147293133
$ws.Cells.Item(3, 10).Value = 0.2032371147293133
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 26.95384733333333
$ws.Cells.Item(3, 14).Value = 80.861542
$ws.Cells.Item(3, 15).Value = 0.6557019380820612
$ws.Cells.Item(3, 16).Value = 0.6557019380820612
$ws.Cells.Item(3, 17).Value = 99.0562335038831
$ws.Cells.Item(3, 18).Value = 891.5061015349479
$ws.Cells.Item(3, 19).Value = 0.133262970018217
$ws.Cells.Item(3, 20).Value = 0.133262970018217

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Inhba"
$ws.Cells.Item(4, 3).Value = "Acvr2a"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 3.675031333333333
$ws.Cells.Item(4, 8).Value = 11.025094
$ws.Cells.Item(4, 9).Value = 0.2032371147293133
$ws.Cells.Item(4, 10).Value = 0.2032371147293133
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 6.212609666666666
$ws.Cells.Item(4, 14).Value = 18.637829
$ws.Cells.Item(4, 15).Value = 0.1511331628692172
$ws.Cells.Item(4, 16).Value = 0.1511331628692172
$ws.Cells.Item(4, 17).Value = 22.83153518676955
$ws.Cells.Item(4, 18).Value = 205.483816680926
$ws.Cells.Item(4, 19).Value = 0.03071586796145508
$ws.Cells.Item(4, 20).Value = 0.03071586796145509

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Inhba"
$ws.Cells.Item(5, 3).Value = "Acvr2a"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 10.108494
$ws.Cells.Item(5, 8).Value = 30.325482
$ws.Cells.Item(5, 9).Value = 0.5590213983169419
$ws.Cells.Item(5, 10).Value = 0.5590213983169419
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 7.940402333333334
$ws.Cells.Item(5, 14).Value = 23.821207
$ws.Cells.Item(5, 15).Value = 0.1931648990487216
$ws.Cells.Item(5, 16).Value = 0.1931648990487216
$ws.Cells.Item(5, 17).Value = 80.26550934408601
$ws.Cells.Item(5, 18).Value = 722.389584096774
$ws.Cells.Item(5, 19).Value = 0.1079833119719673
$ws.Cells.Item(5, 20).Value = 0.1079833119719673

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Inhba"
$ws.Cells.Item(6, 3).Value = "Acvr2a"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 10.108494
$ws.Cells.Item(6, 8).Value = 30.325482
$ws.Cells.Item(6, 9).Value = 0.5590213983169419
$ws.Cells.Item(6, 10).Value = 0.5590213983169419
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 26.95384733333333
$ws.Cells.Item(6, 14).Value = 80.861542
$ws.Cells.Item(6, 15).Value = 0.6557019380820612
$ws.Cells.Item(6, 16).Value = 0.6557019380820612
$ws.Cells.Item(6, 17).Value = 272.462804045916
$ws.Cells.Item(6, 18).Value = 2452.165236413244
$ws.Cells.Item(6, 19).Value = 0.3665514143057627
$ws.Cells.Item(6, 20).Value = 0.3665514143057627

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Inhba"
$ws.Cells.Item(7, 3).Value = "Acvr2a"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 10.108494
$ws.Cells.Item(7, 8).Value = 30.325482
$ws.Cells.Item(7, 9).Value = 0.5590213983169419
$ws.Cells.Item(7, 10).Value = 0.5590213983169419
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 6.212609666666666
$ws.Cells.Item(7, 14).Value = 18.637829
$ws.Cells.Item(7, 15).Value = 0.1511331628692172
$ws.Cells.Item(7, 16).Value = 0.1511331628692172
$ws.Cells.Item(7, 17).Value = 62.800127539842
$ws.Cells.Item(7, 18).Value = 565.201147858578
$ws.Cells.Item(7, 19).Value = 0.0844866720392119
$ws.Cells.Item(7, 20).Value = 0.08448667203921191

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Inhba"
$ws.Cells.Item(8, 3).Value = "Acvr2a"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.298956
$ws.Cells.Item(8, 8).Value = 12.896868
$ws.Cells.Item(8, 9).Value = 0.2377414869537448
$ws.Cells.Item(8, 10).Value = 0.2377414869537448
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 7.940402333333334
$ws.Cells.Item(8, 14).Value = 23.821207
$ws.Cells.Item(8, 15).Value = 0.1931648990487216
$ws.Cells.Item(8, 16).Value = 0.1931648990487216
$ws.Cells.Item(8, 17).Value = 34.13544025329734
$ws.Cells.Item(8, 18).Value = 307.218962279676
$ws.Cells.Item(8, 19).Value = 0.04592331032711308
$ws.Cells.Item(8, 20).Value = 0.04592331032711308

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Inhba"
$ws.Cells.Item(9, 3).Value = "Acvr2a"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.298956
$ws.Cells.Item(9, 8).Value = 12.896868
$ws.Cells.Item(9, 9).Value = 0.2377414869537448
$ws.Cells.Item(9, 10).Value = 0.2377414869537448
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 26.95384733333333
$ws.Cells.Item(9, 14).Value = 80.861542
$ws.Cells.Item(9, 15).Value = 0.6557019380820612
$ws.Cells.Item(9, 16).Value = 0.6557019380820612
$ws.Cells.Item(9, 17).Value = 115.8734037167173
$ws.Cells.Item(9, 18).Value = 1042.860633450456
$ws.Cells.Item(9, 19).Value = 0.1558875537580815
$ws.Cells.Item(9, 20).Value = 0.1558875537580815

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Inhba"
$ws.Cells.Item(10, 3).Value = "Acvr2a"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.298956
$ws.Cells.Item(10, 8).Value = 12.896868
$ws.Cells.Item(10, 9).Value = 0.2377414869537448
$ws.Cells.Item(10, 10).Value = 0.2377414869537448
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 6.212609666666666
$ws.Cells.Item(10, 14).Value = 18.637829
$ws.Cells.Item(10, 15).Value = 0.1511331628692172
$ws.Cells.Item(10, 16).Value = 0.1511331628692172
$ws.Cells.Item(10, 17).Value = 26.70773560217467
$ws.Cells.Item(10, 18).Value = 240.369620419572
$ws.Cells.Item(10, 19).Value = 0.03593062286855018
$ws.Cells.Item(10, 20).Value = 0.03593062286855019
